$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") is stored as plain text (e.g. "30.300.50", "1.001").
# Assigning a plain numeric-looking string via .Value would let Excel
# auto-convert it to a floating point number, losing the original text
# formatting/precision. For each such new price we switch the cell to Text
# format first so it is stored as a string, matching the source data.
# Prices that still contain the thousand-separator dot (e.g. "30.314.79")
# are never parsed as numbers, so they are left with their existing format.

$ws.Range('D2').Value = '30.314.79'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.930.47'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9995'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '249.49'
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.7200'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '28.04'
$ws.Range('E8').Value = '  -2.42%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3210'
$ws.Range('E9').Value = '  -4.09%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07114'
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.7895'
$ws.Range('E11').Value = '  -3.30%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08009'
$ws.Range('E12').Value = '  -1.82%  '
$ws.Range('D13').Value = '1.929.27'
$ws.Range('E13').Value = '  -0.60%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.385'
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '94.87'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.69'
$ws.Range('E16').Value = '  -1.44%  '
$ws.Range('D17').Value = '30.302.76'
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '257.11'
$ws.Range('E18').Value = '  +1.40%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000008098'
$ws.Range('E19').Value = '  -3.26%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.761'
$ws.Range('E20').Value = '  -1.81%  '
$ws.Range('D21').Value = '2.183.74'
$ws.Range('E21').Value = '  -0.38%  '
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.9993'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.832'
$ws.Range('E24').Value = '  -2.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.570'
$ws.Range('E25').Value = '  -2.83%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.73'
$ws.Range('E26').Value = '  +2.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.13'
$ws.Range('E27').Value = '  -1.46%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.290'
$ws.Range('E28').Value = '  -4.75%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.1283'
$ws.Range('E29').Value = '  -2.47%  '
$ws.Range('E30').Value = '  +0.76%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.533'
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.419'
$ws.Range('E32').Value = '  -1.64%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.155'
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05128'
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.283'
$ws.Range('E35').Value = '  +0.95%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7494'
$ws.Range('E36').Value = '  -2.28%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.773'
$ws.Range('E37').Value = '  +0.75%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01989'
$ws.Range('E38').Value = '  -0.60%  '
$ws.Range('E39').Value = '  -1.58%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '78.33'
$ws.Range('E40').Value = '  -3.07%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.395'
$ws.Range('E41').Value = '  -2.56%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4520'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.999'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8456'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.9998'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.41'
$ws.Range('E46').Value = '  -1.78%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.793'
$ws.Range('E47').Value = '  -0.61%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.489'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '36.89'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '962.21'
$ws.Range('E50').Value = '  +8.55%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4217'
$ws.Range('E51').Value = '  +0.71%  '
